# loginTestData.xlsx - record the pass/fail result of each login attempt
# in column D, and leave the selection on D3 (work in progress).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 are the invalid admin1/admin2/admin3 login attempts -> Fail.
$ws.Range("D2").Value = "Fail"
$ws.Range("D3").Value = "Fail"
$ws.Range("D4").Value = "Fail"

# Row 5 is the valid tomsmith login -> Pass.
$ws.Range("D5").Value = "Pass"

# Leave the active selection on D3, matching where editing left off.
[void]$ws.Range("D3").Select()
